$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G10").Value = "x"
$ws.Range("G11").Value = "x"

$ws.Range("G12").Select()
